# Updated policies and graphs
# - Adds 12 new calendar rows (9/30/2020 .. 10/11/2020) to "Converted Data"
# - The "Industries" policy weight (U7) flips 1 -> 0, dropping total weight
#   (X7) from 13 -> 12, which rescales every LockdownEffectiveness value in
#   column X (rows 24-221) by 13/12, and the newly appended rows get the
#   freshly-recomputed weight (1/12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Append 12 new rows, each a clone of the last existing data row (221)
#    -- same flag pattern (only D and F flagged) -- but with its own date
#    label in column A and a placeholder in column X (fixed up in step 3).
# ---------------------------------------------------------------------
$newDates = @("9/30/2020", "10/1/2020", "10/2/2020", "10/3/2020", "10/4/2020", "10/5/2020", "10/6/2020", "10/7/2020", "10/8/2020", "10/9/2020", "10/10/2020", "10/11/2020")

$destRow = 222
foreach ($d in $newDates) {
    $ws.Range("A221:X221").Copy($ws.Cells.Item($destRow, 1))

    # Force the date text into column A as a literal shared string (not an
    # auto-converted date serial) while keeping the bold/border/center style
    # that the Copy() above already brought along.
    $cell = $ws.Cells.Item($destRow, 1)
    $cell.Formula = '=TEXT("' + $d + '","m/d/yyyy")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)

    $destRow = $destRow + 1
}

# ---------------------------------------------------------------------
# 2. Flip the "Industries" weight off and drop the total weight.
# ---------------------------------------------------------------------
$ws.Range("U7").Value = 0
$ws.Range("X7").Value = 12

# ---------------------------------------------------------------------
# 3. Recompute column X (LockdownEffectiveness) for every data row now
#    that the total weight denominator changed from 13 to 12. The numeric
#    numerator (sum of flagged weights) is unchanged, so old*13/12 ==
#    new value for every pre-existing row; the freshly appended rows
#    (same flag pattern as row 221) land on 1/12 directly.
# ---------------------------------------------------------------------
for ($r = 24; $r -le 221; $r++) {
    $old = $ws.Cells.Item($r, 24).Value()
    $ws.Cells.Item($r, 24).Value = $old * 13 / 12
}

for ($r = 222; $r -le 233; $r++) {
    $ws.Cells.Item($r, 24).Value = 1 / 12
}
